$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting A:E -> B:F
$ws.Columns("A:A").Insert()

# Match the header formatting (bold, bordered, centered) used by the rest
# of row 1 so the new "ID" header looks consistent with its neighbors.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# New header for the inserted ID column
$ws.Range("A1").Value = "ID"

# ID values for the newly inserted column A (rows 2-25)
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
